$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "SCD0026"

# Update the TC_ID value in B2
$ws.Range("B2").Value = "SCD0026-002"

# Update the selected cell to B3
$ws.Range("B3").Select()
